# Provision OA LEGAL - FR: update salutation line to split gender/client,
# and merge the parenthesised TVA fragment into a single run.

$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# "Mr/Ms {client}" -> "{gender}" + " {client}" (two runs, same formatting)
$rng = $d.Content
$found1 = $rng.Find.Execute("Mr/Ms {client}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng.Text = "{gender}"
    $rng.Collapse(0)
    $rng.InsertAfter(" {client}")
    # Force the two adjacent, identically-formatted runs to stay distinct
    # (toggling a character property and reverting it breaks the auto-merge).
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# --- Change 2 --------------------------------------------------------
# "(" + "{TVA}" + "), à titre ... suivantes :" (3 runs) -> 1 merged run
# Locate the span with anchor searches (avoids hard-coding the nbsp before
# the final colon) and collapse it into a single run via Find/Replace.
$rngTva = $d.Content
[void]$rngTva.Find.Execute("{TVA}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tvaStart = $rngTva.Start

$rngTail = $d.Content
[void]$rngTail.Find.Execute("suivantes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailEnd = $rngTail.End + 2

$full = $d.Range($tvaStart - 1, $tailEnd)
$fullText = $full.Text
$found2 = $full.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, $fullText, 2)
